$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value of B2 with the new "modif 14h10" text
$ws.Range("B2").Value = "Donnée B2 - modif 14h10"

# Reflect the active selection recorded in the saved file (user was editing B2)
$ws.Range("B2").Select()
